# TestCasesRobustness.xlsx — "Edit new date and expect result"
#
# The ActualRank test input (column B, VisitsPerMonth) was bumped from 4
# visits/month to 15 (and a few outlier rows to 29/30/31) and the
# expected-rank / expected-result columns were refreshed to match the
# membership-ranking rules for the new data. Row 20's Spending value was
# also fixed from 5000 to 50000 (with its expected columns updated), and
# the active selection moved to J15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (VisitsPerMonth) bulk update: 4 -> 15 for most rows ---
$bRows15 = 2,3,4,5,6,7,8,15,16,17,18,19,20
foreach ($r in $bRows15) {
    $ws.Cells.Item($r, 2).Value = 15
}

# Rows 3-8 previously carried a "wrap text / vertical center" style (s=2);
# after the edit they fall back to the sheet's default (unstyled) cell.
$ws.Range("B3:B8").ClearFormats()

# A few rows got distinct new VisitsPerMonth outliers instead of 15.
$ws.Cells.Item(12, 2).Value = 29
$ws.Cells.Item(13, 2).Value = 30
$ws.Cells.Item(14, 2).Value = 31

# --- ExpectedRank (column D) refresh: rows that used to expect "Gold"
#     now expect "Standard" for the updated VisitsPerMonth value ---
$ws.Range("D5").Value = "Standard"
$ws.Range("D6").Value = "Standard"
$ws.Range("D7").Value = "Standard"
$ws.Range("D18").Value = "Standard"
$ws.Range("D19").Value = "Standard"

# --- Row 20: Spending input fixed, and its expected/actual columns
#     refreshed accordingly ---
$ws.Range("A20").Value = 50000
$ws.Range("E20").Value = "Silver"
$ws.Range("F20").Value = "FAIL"

# --- Move the active selection to J15 ---
$ws.Range("J15").Select()
